# Re-serializes word/document.xml and word/numbering.xml so that every
# element's attributes are emitted in alphabetical order (namespace
# declarations first - sorted among themselves - followed by the rest of
# the attributes, also sorted), matching the output of the upgraded
# Apache POI (XMLBeans) writer used after the POI 3.15 upgrade. This is
# a pure re-ordering of attributes: no element, attribute, or text value
# is added, removed, or changed.

function Sort-TagAttrs($sta_tagInner) {
    # $sta_tagInner is the raw text between '<' and '>' (exclusive), e.g.
    #   w:document xmlns:a="..." xmlns:b="..."
    #   w:r/
    #   /w:r
    #   ?xml version="1.0" encoding="UTF-8"?
    if ($sta_tagInner.StartsWith("/")) {
        return $sta_tagInner
    }
    if ($sta_tagInner.StartsWith("?") -or $sta_tagInner.StartsWith("!")) {
        return $sta_tagInner
    }

    $sta_selfClose = $false
    $sta_work = $sta_tagInner
    if ($sta_work.EndsWith("/")) {
        $sta_selfClose = $true
        $sta_work = $sta_work.Substring(0, $sta_work.Length - 1)
    }

    $sta_m = [regex]::Match($sta_work, '^(\S+)(\s*)([\s\S]*)$')
    $sta_name = $sta_m.Groups[1].Value
    $sta_rest = $sta_m.Groups[3].Value

    if ($sta_rest.Trim().Length -eq 0) {
        $sta_out = $sta_name
        if ($sta_selfClose) { $sta_out += "/" }
        return $sta_out
    }

    $sta_attrMatches = [regex]::Matches($sta_rest, '([\w:.\-]+)=("[^"]*"|''[^'']*'')')
    $sta_nsAttrs = @()
    $sta_otherAttrs = @()
    foreach ($sta_am in $sta_attrMatches) {
        $sta_k = $sta_am.Groups[1].Value
        $sta_v = $sta_am.Groups[2].Value
        $sta_entry = @{ K = $sta_k; Text = "$sta_k=$sta_v" }
        if ($sta_k -eq "xmlns" -or $sta_k.StartsWith("xmlns:")) {
            $sta_nsAttrs += $sta_entry
        } else {
            $sta_otherAttrs += $sta_entry
        }
    }

    $sta_nsSorted = $sta_nsAttrs | Sort-Object { $_.K }
    $sta_otherSorted = $sta_otherAttrs | Sort-Object { $_.K }

    $sta_parts = @()
    foreach ($sta_e in $sta_nsSorted) { $sta_parts += $sta_e.Text }
    foreach ($sta_e in $sta_otherSorted) { $sta_parts += $sta_e.Text }

    $sta_out = $sta_name
    if ($sta_parts.Length -gt 0) {
        $sta_out += " " + ($sta_parts -join " ")
    }
    if ($sta_selfClose) { $sta_out += "/" }
    return $sta_out
}

function Sort-XmlTagAttrsInString($sxt_xmlText) {
    $sxt_sb = New-Object System.Text.StringBuilder
    $sxt_i = 0
    $sxt_n = $sxt_xmlText.Length
    while ($sxt_i -lt $sxt_n) {
        $sxt_ch = $sxt_xmlText[$sxt_i]
        if ($sxt_ch -eq '<') {
            $sxt_k = $sxt_i + 1
            $sxt_inQuote = ""
            while ($sxt_k -lt $sxt_n) {
                $sxt_c = $sxt_xmlText[$sxt_k]
                if ($sxt_inQuote -ne "") {
                    if ($sxt_c -eq $sxt_inQuote) { $sxt_inQuote = "" }
                } elseif ($sxt_c -eq '"' -or $sxt_c -eq "'") {
                    $sxt_inQuote = $sxt_c
                } elseif ($sxt_c -eq '>') {
                    break
                }
                $sxt_k++
            }
            $sxt_tagInner = $sxt_xmlText.Substring($sxt_i + 1, $sxt_k - $sxt_i - 1)
            $sxt_newTagInner = Sort-TagAttrs $sxt_tagInner
            [void]$sxt_sb.Append("<")
            [void]$sxt_sb.Append($sxt_newTagInner)
            [void]$sxt_sb.Append(">")
            $sxt_i = $sxt_k + 1
        } else {
            $sxt_nxt = $sxt_xmlText.IndexOf('<', $sxt_i)
            if ($sxt_nxt -lt 0) {
                [void]$sxt_sb.Append($sxt_xmlText.Substring($sxt_i))
                $sxt_i = $sxt_n
            } else {
                [void]$sxt_sb.Append($sxt_xmlText.Substring($sxt_i, $sxt_nxt - $sxt_i))
                $sxt_i = $sxt_nxt
            }
        }
    }
    return $sxt_sb.ToString()
}

function Reorder-PackagePart($rpp_fullXml, $rpp_partName) {
    $rpp_marker = 'pkg:name="' + $rpp_partName + '"'
    $rpp_partIdx = $rpp_fullXml.IndexOf($rpp_marker)
    if ($rpp_partIdx -lt 0) {
        return $rpp_fullXml
    }
    $rpp_dataOpenTag = "<pkg:xmlData>"
    $rpp_dataStart = $rpp_fullXml.IndexOf($rpp_dataOpenTag, $rpp_partIdx) + $rpp_dataOpenTag.Length
    $rpp_dataCloseTag = "</pkg:xmlData>"
    $rpp_dataEnd = $rpp_fullXml.IndexOf($rpp_dataCloseTag, $rpp_dataStart)

    $rpp_before = $rpp_fullXml.Substring(0, $rpp_dataStart)
    $rpp_inner = $rpp_fullXml.Substring($rpp_dataStart, $rpp_dataEnd - $rpp_dataStart)
    $rpp_after = $rpp_fullXml.Substring($rpp_dataEnd)

    $rpp_newInner = Sort-XmlTagAttrsInString $rpp_inner

    return $rpp_before + $rpp_newInner + $rpp_after
}

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$xml = Reorder-PackagePart $xml "/word/document.xml"
$xml = Reorder-PackagePart $xml "/word/numbering.xml"

$d.WordOpenXML = $xml
